$d = $word.ActiveDocument

$replacements = @(
    @{old = "276×8=2208"; new = "167×6=1002"},
    @{old = "977×5=4885"; new = "155×8=1240"},
    @{old = "486×7=3402"; new = "549×9=4941"},
    @{old = "768×8=6144"; new = "978×3=2934"},
    @{old = "310×8=2480"; new = "783×9=7047"},
    @{old = "413×9=3717"; new = "374×8=2992"},
    @{old = "565×2=1130"; new = "703×2=1406"},
    @{old = "710×8=5680"; new = "556×7=3892"},
    @{old = "275×7=1925"; new = "648×4=2592"},
    @{old = "688×7=4816"; new = "638×9=5742"},
    @{old = "661×5=3305"; new = "207×6=1242"},
    @{old = "147×5=735";  new = "978×7=6846"},
    @{old = "725×7=5075"; new = "143×2=286"},
    @{old = "930×2=1860"; new = "392×8=3136"},
    @{old = "206×9=1854"; new = "605×9=5445"},
    @{old = "959×9=8631"; new = "749×8=5992"},
    @{old = "138×8=1104"; new = "548×4=2192"},
    @{old = "655×3=1965"; new = "731×3=2193"},
    @{old = "354×7=2478"; new = "782×9=7038"},
    @{old = "114×5=570";  new = "607×8=4856"},
    @{old = "611×8=4888"; new = "682×7=4774"},
    @{old = "249×7=1743"; new = "635×9=5715"},
    @{old = "318×5=1590"; new = "378×5=1890"},
    @{old = "926×3=2778"; new = "298×8=2384"},
    @{old = "383×3=1149"; new = "720×4=2880"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
